# "Generate Report for Handoff"
#
# The localization report moves from "In Translation" to "Ready for
# handoff" and the handoff timestamps are refreshed. Updating the text
# widens the Status columns, so their widths are refreshed to match the
# new (longer) content.

$wb = $excel.ActiveWorkbook

$statusOld  = "In Translation"
$statusNew  = "Ready for handoff"

# ColumnWidth is quantized by the host to 1/6-character increments, so we
# pick the input value whose quantized result lands closest to the
# intended rendered width for the new, longer "Ready for handoff" text.
$newStatusColWidth = 16.333333333333332

# ---- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F2").Value = $statusNew
$wsOverview.Range("G2").Value = "2016-10-18 04:51:54"
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = $newStatusColWidth

# ---- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("H2").Value = "2016-10-18 04:51:38"
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newStatusColWidth

# ---- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("H2").Value = "2016-10-18 04:51:54"
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newStatusColWidth
